$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "7-Day Avg Hours" column (B) was regenerated from a new random data
# stream in the generating code, so every sliding-average value in B2:B91
# changes. Apply the new values exactly as produced by the regenerated
# analysis (B92 is unaffected/not covered by this change).
$newAverages = @(
    4.65, 5.487500000000001, 5.752727272727273, 5.435833333333334, 5.089375, 5.253157894736842, 5.4736363636363645, 5.464615384615385,
    5.29913043478261, 5.034800000000001, 5.490625, 5.8032258064516125, 5.68470588235294, 5.446969696969696, 5.5216666666666665, 5.666363636363637,
    5.60027027027027, 5.32764705882353, 5.2592105263157904, 5.355526315789475, 5.535384615384618, 5.460526315789475, 5.19810810810811, 5.27451612903226,
    5.16, 5.158125000000001, 5.003571428571428, 5.131034482758621, 5.163437500000001, 5.4175, 5.579189189189188, 5.551142857142857,
    5.693714285714284, 5.644, 5.344318181818181, 5.3252500000000005, 5.087894736842105, 4.973030303030302, 5.264411764705883, 4.796774193548389,
    4.748928571428571, 4.558333333333333, 4.5495833333333335, 5.261153846153845, 5.253571428571428, 5.126956521739131, 5.362727272727273, 5.558571428571429,
    6.17857142857143, 6.030833333333334, 5.531538461538462, 5.3852, 5.340370370370372, 5.696896551724137, 5.354375, 4.913749999999999,
    5.272571428571428, 5.1270588235294134, 4.954102564102565, 4.855581395348839, 4.560250000000002, 4.86075, 5.052999999999999, 4.827499999999998,
    4.928333333333333, 5.1930000000000005, 5.5592, 5.841724137931035, 5.478666666666666, 5.281818181818183, 5.310882352941177, 5.4534375,
    5.151875000000001, 4.793428571428573, 4.620588235294117, 4.835312499999998, 5.566333333333332, 5.586428571428572, 5.846129032258064, 6.0687878787878775,
    6.4734375, 6.401333333333334, 6.41551724137931, 5.945769230769231, 5.652068965517241, 5.177777777777778, 4.966666666666667, 4.5484,
    4.6328, 4.489166666666667
)

for ($i = 0; $i -lt $newAverages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newAverages[$i]
}

# The bestFit column widths were recalculated by the authoring tool for the
# new cell contents (col A: 11.3046875 -> 11.69921875, col B: 15.83984375 ->
# 16.53125). Apply the closest widths obtainable through the ColumnWidth
# property (this engine quantizes stored width to whole-pixel steps).
$ws.Columns.Item(1).ColumnWidth = 10.865885416666668
$ws.Columns.Item(2).ColumnWidth = 15.697916666666668
